$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlPasteFormats = -4122

function Copy-Format($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats) | Out-Null
}

# ===========================================================================
# Ch 08 06 "Understanding advanced time intelligence calculations"
# ===========================================================================

# ---------------------------------------------------------------------------
# Row 37: this row already existed (E/F/G/H/I filled in) - add the B/C/D
# cells that complete it (Key formula, Chapter number, Section name)
# ---------------------------------------------------------------------------
$ws.Range("C37").Value = 8
$ws.Range("D37").Value = "Time intelligence calculations"

# ---------------------------------------------------------------------------
# Row 38: same - add B/C/D, plus the new I38 function reference
# ---------------------------------------------------------------------------
$ws.Range("C38").Value = 8
$ws.Range("D38").Value = "Time intelligence calculations"
$ws.Range("I38").Value = "LASTNONBLANK, FIRSTNONBLANK, PARALLELPERIOD"

# B37:B38 share the same "Key" formula - assigning the formula across both
# cells at once lets them be written out as a shared formula group.
$ws.Range("B37:B38").Formula = '=_xlfn.CONCAT(TEXT(C37,"00"),TEXT(E37,"00"),TEXT(G37,"00"))'

# ---------------------------------------------------------------------------
# Row 39 (brand new row) - new subsection "Understanding periods to date"
# ---------------------------------------------------------------------------
$ws.Range("C39").Value = 8
$ws.Range("D39").Value = "Time intelligence calculations"
$ws.Range("E39").Value = 6
$ws.Range("F39").Value = "Understanding advanced time intelligence calculations"
$ws.Range("G39").Value = 1
$ws.Range("H39").Value = "Understanding periods to date"
$ws.Range("I39").Value = "DATESYTD"

# ---------------------------------------------------------------------------
# Row 40 (brand new row) - "Understanding DATEADD"
# ---------------------------------------------------------------------------
$ws.Range("C40").Value = 8
$ws.Range("D40").Value = "Time intelligence calculations"
$ws.Range("E40").Value = 6
$ws.Range("F40").Value = "Understanding advanced time intelligence calculations"
$ws.Range("G40").Value = 2
$ws.Range("H40").Value = "Understanding DATEADD"

# ---------------------------------------------------------------------------
# Row 41 (brand new row) - "Understanding FIRSTDATE, LASTDATE, ..."
# ---------------------------------------------------------------------------
$ws.Range("C41").Value = 8
$ws.Range("D41").Value = "Time intelligence calculations"
$ws.Range("E41").Value = 6
$ws.Range("F41").Value = "Understanding advanced time intelligence calculations"
$ws.Range("G41").Value = 3
$ws.Range("H41").Value = "Understanding FIRSTDATE, LASTDATE, FIRSTNONBLANK, LASTNONBLANK"
$ws.Range("I41").Value = "FIRSTDATE, LASTDATE, FIRSTNONBLANK, LASTNONBLANK"

# B39:B41 share the same "Key" formula - write it as one shared formula group.
$ws.Range("B39:B41").Formula = '=_xlfn.CONCAT(TEXT(C39,"00"),TEXT(E39,"00"),TEXT(G39,"00"))'

# ---------------------------------------------------------------------------
# Row 42 (brand new row) - "Using drillthrough with time intelligence"
# ---------------------------------------------------------------------------
$ws.Range("C42").Value = 8
$ws.Range("D42").Value = "Time intelligence calculations"
$ws.Range("E42").Value = 6
$ws.Range("F42").Value = "Understanding advanced time intelligence calculations"
$ws.Range("G42").Value = 4
$ws.Range("H42").Value = "Using drillthrough with time intelligence"
$ws.Range("B42").Formula = '=_xlfn.CONCAT(TEXT(C42,"00"),TEXT(E42,"00"),TEXT(G42,"00"))'

# ---------------------------------------------------------------------------
# Copy cell formatting (borders / number formats) from row 36, which is the
# last fully-styled row of this block, down across the new/extended rows.
# ---------------------------------------------------------------------------
Copy-Format "B36" "B37:B42"
Copy-Format "C36" "C37:C42"
Copy-Format "E36" "E37:E42"
Copy-Format "G36" "G37:G42"
Copy-Format "I36" "I38:I39"
Copy-Format "I36" "I41"

# ---------------------------------------------------------------------------
# Column widths widened (bestFit) to accommodate the new, longer strings in
# columns F (50.7109375 chars) and H (66.42578125 chars)
# ---------------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 49.833333333333336
$ws.Columns.Item(8).ColumnWidth = 65.66666666666667

# ---------------------------------------------------------------------------
# Sheet view: scroll position / selection change recorded in the diff
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("H49").Select() | Out-Null
